$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data cells whose values changed ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $false

$ws.Range("D13").Value = 8

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# Two new "line" entries were inserted into the name list right after
# "line6" (and before the "extr*" entries). Because of that, rows 8-15
# (which used to refer to extr1..extr8) now line up with line7, line8,
# extr1..extr6, and two brand-new rows are appended for extr7 / extr8.
$ws.Range("B8").Value = "line7"
$ws.Range("B9").Value = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"

# --- Add two new rows (16 and 17) for extr7 / extr8 ---
# Copy formatting from the existing "A" column cell (bold, centered, bordered)
# so the new rows visually match the rest of the table.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
